$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying style from E1 (bold header style) so it matches s="1"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (rows 2-98)
$timeTaken = @(
    "2021-10-05 10:52:33.773909",
    "2021-10-05 10:52:33.773921",
    "2021-10-05 10:52:33.773924",
    "2021-10-05 10:52:33.773927",
    "2021-10-05 10:52:33.773930",
    "2021-10-05 10:52:33.773932",
    "2021-10-05 10:52:33.773935",
    "2021-10-05 10:52:33.773938",
    "2021-10-05 10:52:33.773941",
    "2021-10-05 10:52:33.773944",
    "2021-10-05 10:52:33.773947",
    "2021-10-05 10:52:33.773949",
    "2021-10-05 10:52:33.773952",
    "2021-10-05 10:52:33.773954",
    "2021-10-05 10:52:33.773957",
    "2021-10-05 10:52:33.773960",
    "2021-10-05 10:52:33.773962",
    "2021-10-05 10:52:33.773965",
    "2021-10-05 10:52:33.773967",
    "2021-10-05 10:52:33.773970",
    "2021-10-05 10:52:33.773972",
    "2021-10-05 10:52:33.773975",
    "2021-10-05 10:52:33.773977",
    "2021-10-05 10:52:33.773980",
    "2021-10-05 10:52:33.773983",
    "2021-10-05 10:52:33.773986",
    "2021-10-05 10:52:33.773988",
    "2021-10-05 10:52:33.773991",
    "2021-10-05 10:52:33.773993",
    "2021-10-05 10:52:33.773996",
    "2021-10-05 10:52:33.773998",
    "2021-10-05 10:52:33.774001",
    "2021-10-05 10:52:33.774004",
    "2021-10-05 10:52:33.774007",
    "2021-10-05 10:52:33.774009",
    "2021-10-05 10:52:33.774012",
    "2021-10-05 10:52:33.774014",
    "2021-10-05 10:52:33.774017",
    "2021-10-05 10:52:33.774019",
    "2021-10-05 10:52:33.774022",
    "2021-10-05 10:52:33.774025",
    "2021-10-05 10:52:33.774028",
    "2021-10-05 10:52:33.774030",
    "2021-10-05 10:52:33.774033",
    "2021-10-05 10:52:33.774036",
    "2021-10-05 10:52:33.774038",
    "2021-10-05 10:52:33.774041",
    "2021-10-05 10:52:33.774043",
    "2021-10-05 10:52:33.774046",
    "2021-10-05 10:52:33.774049",
    "2021-10-05 10:52:33.774051",
    "2021-10-05 10:52:33.774054",
    "2021-10-05 10:52:33.774057",
    "2021-10-05 10:52:33.774060",
    "2021-10-05 10:52:33.774062",
    "2021-10-05 10:52:33.774065",
    "2021-10-05 10:52:33.774067",
    "2021-10-05 10:52:33.774070",
    "2021-10-05 10:52:33.774072",
    "2021-10-05 10:52:33.774075",
    "2021-10-05 10:52:33.774077",
    "2021-10-05 10:52:33.774080",
    "2021-10-05 10:52:33.774082",
    "2021-10-05 10:52:33.774085",
    "2021-10-05 10:52:33.774089",
    "2021-10-05 10:52:33.774091",
    "2021-10-05 10:52:33.774094",
    "2021-10-05 10:52:33.774097",
    "2021-10-05 10:52:33.774100",
    "2021-10-05 10:52:33.774102",
    "2021-10-05 10:52:33.774104",
    "2021-10-05 10:52:33.774107",
    "2021-10-05 10:52:33.774110",
    "2021-10-05 10:52:33.774112",
    "2021-10-05 10:52:33.774115",
    "2021-10-05 10:52:33.774118",
    "2021-10-05 10:52:33.774122",
    "2021-10-05 10:52:33.774126",
    "2021-10-05 10:52:33.774128",
    "2021-10-05 10:52:33.774131",
    "2021-10-05 10:52:33.774134",
    "2021-10-05 10:52:33.774137",
    "2021-10-05 10:52:33.774139",
    "2021-10-05 10:52:33.774142",
    "2021-10-05 10:52:33.774145",
    "2021-10-05 10:52:33.774148",
    "2021-10-05 10:52:33.774150",
    "2021-10-05 10:52:33.774153",
    "2021-10-05 10:52:33.774156",
    "2021-10-05 10:52:33.774159",
    "2021-10-05 10:52:33.774162",
    "2021-10-05 10:52:33.774164",
    "2021-10-05 10:52:33.774168",
    "2021-10-05 10:52:33.774172",
    "2021-10-05 10:52:33.774174",
    "2021-10-05 10:52:33.774177",
    "2021-10-05 10:52:33.774180"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

